# Adapt column header formatting to respective input file names.
# "_old" -> "_FV2304" (left/"before" block, columns A:J)
# "_new" -> "_FV2310" (right/"after" block, columns L:U)
# Also freezes the header row and wraps the data range in a native Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

$leftCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Range($leftCols[$i] + "1").Value  = $fields[$i] + "_FV2304"
    $ws.Range($rightCols[$i] + "1").Value = $fields[$i] + "_FV2310"
}

# K1 ("diff") is unchanged.

# Freeze the header row (Excel's "Freeze Top Row").
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

# Turn the used range into a native table ("Table1"), matching the exported
# AHB-diff table (21 columns, A1:U57, with autofilter).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"

Write-Output "done"
